# Refresh the crypto price/volume table in-place, cell by cell, to
# match the latest values from the data source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.261.09'
$ws.Range('E2').Value = '  +0.34%  '

$ws.Range('D3').Value = '1.896.07'
$ws.Range('E3').Value = '  +2.17%  '

$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').Value = '''244.81'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E5').Value = '  +2.80%  '

$ws.Range('D6').Value = '''0.655'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E6').Value = '  +5.43%  '

$ws.Range('E7').Value = '  -0.28%  '

$ws.Range('D8').Value = '''41.50'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E8').Value = '  -1.45%  '

$ws.Range('E9').Value = '  +7.11%  '

$ws.Range('D10').Value = '''52.21'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E10').Value = '  +11.94%  '

$ws.Range('D11').Value = '''0.0713'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E11').Value = '  +2.92%  '

$ws.Range('D12').Value = '''0.0995'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E12').Value = '  +0.68%  '

$ws.Range('D13').Value = '2.170.06'
$ws.Range('E13').Value = '  +2.21%  '

$ws.Range('D14').Value = '''12.04'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E14').Value = '  +5.21%  '

$ws.Range('D15').Value = '''0.697'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E15').Value = '  +3.12%  '

$ws.Range('D16').Value = '1.890.57'
$ws.Range('E16').Value = '  +1.70%  '

$ws.Range('E17').Value = '  +2.72%  '

$ws.Range('D18').Value = '35.260.03'
$ws.Range('E18').Value = '  +0.48%  '

$ws.Range('D19').Value = '''71.39'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E19').Value = '  +2.14%  '

$ws.Range('D20').Value = '0.0₃0819'
$ws.Range('E20').Value = '  +3.25%  '

$ws.Range('D21').Value = '''240.26'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E21').Value = '  -0.24%  '

$ws.Range('D22').Value = '''12.48'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E22').Value = '  +2.11%  '

$ws.Range('D23').Value = '''4.79'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E23').Value = '  +1.58%  '

$ws.Range('E24').Value = '  -0.30%  '

$ws.Range('D25').Value = '''2.46'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E25').Value = '  +33.24%  '

$ws.Range('D26').Value = '''2.30'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E26').Value = '  +1.06%  '

$ws.Range('D27').Value = '''170.53'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E27').Value = '  +1.54%  '

$ws.Range('D28').Value = '''8.47'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E28').Value = '  +6.48%  '

$ws.Range('D29').Value = '''18.29'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E29').Value = '  +3.67%  '

$ws.Range('E30').Value = '  +2.27%  '

$ws.Range('D31').Value = '''4.13'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E31').Value = '  +3.81%  '

$ws.Range('D32').Value = '''0.0563'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E32').Value = '  +1.50%  '

$ws.Range('D33').Value = '''0.939'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E33').Value = '  +10.66%  '

$ws.Range('E35').Value = '  +2.96%  '

$ws.Range('E36').Value = '  -3.52%  '

$ws.Range('E37').Value = '  +0.03%  '

$ws.Range('E38').Value = '  +1.59%  '

$ws.Range('D39').Value = '''0.0210'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E39').Value = '  +4.47%  '

$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '''1.09'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E40').Value = '  +1.83%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.0649'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E41').Value = '  +16.67%  '

$ws.Range('D42').Value = '''16.30'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E42').Value = '  +9.88%  '

$ws.Range('D43').Value = '''89.84'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E43').Value = '  -0.57%  '

$ws.Range('D44').Value = '1.340.19'
$ws.Range('E44').Value = '  -0.12%  '

$ws.Range('E45').Value = '  +2.37%  '

$ws.Range('D46').Value = '''47.58'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E46').Value = '  +37.08%  '

$ws.Range('E47').Value = '  +1.46%  '

$ws.Range('E48').Value = '  -0.66%  '

$ws.Range('D49').Value = '''6.55'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E49').Value = '  +1.05%  '

$ws.Range('D50').Value = '2.080.35'
$ws.Range('E50').Value = '  +2.06%  '

$ws.Range('D51').Value = '''11.22'  # forced text (leading apostrophe), matches original text-typed price cell
$ws.Range('E51').Value = '  -10.51%  '
